$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "PerciseHit"
$ws.Range("B3").Value = "Heal"
$ws.Range("B6").Value = "Bats"
$ws.Range("B8").Value = "DoubleAttack"
$ws.Range("B9").Value = "FireBall"
$ws.Range("B7").Value = "ShieldBash"
$ws.Range("B5").Value = "PoisonAttack"
$ws.Range("B4").Value = "DrainLife"
$ws.Range("E5").Value = "BodySlam"

$ws.Range("E5").Select()
